$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer-date text (8/30/2018 ->
#    9/4/2018) on the slide master and every slide layout. The placeholder
#    is identified dynamically via PlaceholderFormat.Type (ppPlaceholderDate
#    = 16) so the script does not depend on a fixed shape index.
# ---------------------------------------------------------------------------
$ppPlaceholderDate = 16
$newDate = "9/4/2018"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isPlaceholder = $false
            try { $isPlaceholder = ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) } catch { $isPlaceholder = $false }
            if ($isPlaceholder) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DateShape $p.SlideMaster.Shapes

for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 2 - Agenda placeholder body text updates.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$agenda = $slide2.Shapes.Item(2)
$agendaText = $agenda.TextFrame.TextRange

# "Early Release" -> "Modified Day"
$agendaText.Paragraphs(2, 1).Text = "Modified Day"

# "Tomorrow: More Simon" -> three runs: "Tomorrow" / ": " / "Simon Game"
$agendaText.Paragraphs(7, 1).Text = "Tomorrow: Simon Game"
$para7 = $agendaText.Paragraphs(7, 1)
$run2 = $para7.Characters(9, 2)
$run3 = $para7.Characters(11, 10)
# Re-assign each segment to itself so the engine materialises it as its own
# run (matching the 3-run split in the target markup) without altering the
# visible formatting.
$run2.Text = $run2.Text
$run3.Text = $run3.Text
